# Zhihong Deng handover sheet -> verification pass by Muhammad Arslan.
# See commit message: "Completed Verification of Zhihing Inheritance Work"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) The continuing member handing over was "Zhihong Deng" (column I, rows
#    4-12, all sharing one string). The new reviewer taking over / verifying
#    the handover is "Muhammad Arslan" - update every occurrence together so
#    the shared string itself is renamed rather than forked.
# ---------------------------------------------------------------------------
$ws.Range("I4:I12").Value = "Muhammad Arslan"

# ---------------------------------------------------------------------------
# 2) Fill in the reviewer-facing columns M (Verified), N (Status/Degree of
#    Inheritance) and O (Comments) for every task row, which were previously
#    left blank.
# ---------------------------------------------------------------------------

# Rows 4-5: font-skeleton tasks - Arslan says he was told about the task but
# has not personally run it yet.
$ws.Range("M4:M5").Value = "Yes, Zhihong told me about this task"
$ws.Range("N4:N5").Value = "Yeah, he highlighted these things to me in repo and gave me a bit overview but I haven't ran it on my system yet, but some modules did worked on my laptop"
$ws.Range("O4:O5").Value = "Can explain it to the new members next semester and we can try to improve it and make it more better in working."

# Rows 6-12: remaining tasks - Arslan understood these parts but still wants
# to test them further himself before handing them on.
$ws.Range("M6:M12").Value = "Yes, Zhihong told me about this task and I understood this part as well"
$ws.Range("N6:N12").Value = "Yes, I have understood this part and but need to further test it on my system to get better idea of it"
$ws.Range("O6:O12").Value = "Will need to put some effort in to digest this work and will try to incorporate all the knowledge base I have regarding this issue it to the new members as well"

# Match the wrapped-text styling already used for the rest of the table.
$ws.Range("M4:O12").WrapText = $true

# ---------------------------------------------------------------------------
# 3) Row heights grow now that M/N/O actually contain wrapped commentary.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 138.75
$ws.Rows.Item(6).RowHeight = 138.75
$ws.Rows.Item(7).RowHeight = 126
$ws.Rows.Item(8).RowHeight = 120

# ---------------------------------------------------------------------------
# 4) Column M ("Verified") now needs to be wide enough for its new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 34.28515625

# Land the cursor/selection where Arslan's review ended up.
$null = $ws.Range("N18").Select()
